$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (age, sex, edu) for each of the 50 subjects
$data = @(
    @(57, "Male", 6),
    @(60, "Male", 6),
    @(61, "Female", 18),
    @(60, "Female", 13),
    @(67, "Female", 5),
    @(76, "Male", 10),
    @(78, "Female", 6),
    @(73, "Male", 20),
    @(59, "Male", 5),
    @(61, "Female", 8),
    @(52, "Male", 5),
    @(53, "Male", 15),
    @(51, "Female", 10),
    @(76, "Male", 8),
    @(50, "Male", 5),
    @(61, "Female", 5),
    @(68, "Male", 8),
    @(57, "Female", 6),
    @(74, "Female", 5),
    @(73, "Female", 6),
    @(63, "Male", 7),
    @(71, "Female", 8),
    @(71, "Female", 10),
    @(61, "Female", 20),
    @(63, "Male", 9),
    @(69, "Female", 16),
    @(78, "Male", 15),
    @(75, "Male", 7),
    @(66, "Male", 7),
    @(65, "Female", 6),
    @(59, "Female", 6),
    @(60, "Male", 5),
    @(63, "Female", 7),
    @(58, "Male", 9),
    @(78, "Female", 13),
    @(77, "Male", 9),
    @(72, "Male", 15),
    @(63, "Male", 7),
    @(56, "Female", 5),
    @(59, "Male", 16),
    @(52, "Male", 5),
    @(79, "Male", 7),
    @(60, "Female", 18),
    @(57, "Female", 11),
    @(54, "Female", 7),
    @(72, "Male", 9),
    @(70, "Male", 8),
    @(54, "Male", 7),
    @(56, "Male", 7),
    @(60, "Female", 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

